# Weekly crime-stat data refresh: CompStat_1 sheet now covers the week of
# 2/13/2023-2/19/2023 (Volume 30, Number 7) instead of 2/6/2023-2/12/2023.
# All Week-to-Date / 28-Day / Year-to-Date figures (and their %-change
# columns) for rows 14-30 are updated with the newly collected counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header (rich-text shared strings) ---
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# --- Crime-complaint data table (rows 14-30, columns C:N) ---
# Row 14: Murder
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = -50
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 11
$ws.Range("H14").Value = -45.454545454545
$ws.Range("I14").Value = 14
$ws.Range("J14").Value = 12
$ws.Range("K14").Value = 16.666666666666
$ws.Range("L14").Value = 180
$ws.Range("M14").Value = 16.666666666666
$ws.Range("N14").Value = -63.157894736842

# Row 15: Rape
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 17
$ws.Range("G15").Value = 18
$ws.Range("H15").Value = -5.555555555555
$ws.Range("I15").Value = 33
$ws.Range("J15").Value = 29
$ws.Range("K15").Value = 13.793103448275
$ws.Range("L15").Value = 32
$ws.Range("M15").Value = 120
$ws.Range("N15").Value = -52.173913043478

# Row 16: Robbery
$ws.Range("C16").Value = 34
$ws.Range("D16").Value = 57
$ws.Range("E16").Value = -40.350877192982
$ws.Range("F16").Value = 144
$ws.Range("G16").Value = 167
$ws.Range("H16").Value = -13.772455089820
$ws.Range("I16").Value = 252
$ws.Range("J16").Value = 283
$ws.Range("K16").Value = -10.954063604240
$ws.Range("L16").Value = 56.521739130434
$ws.Range("M16").Value = -32.978723404255
$ws.Range("N16").Value = -87.659157688540

# Row 17: Fel. Assault
$ws.Range("C17").Value = 65
$ws.Range("D17").Value = 73
$ws.Range("E17").Value = -10.958904109589
$ws.Range("G17").Value = 266
$ws.Range("H17").Value = -14.661654135338
$ws.Range("I17").Value = 452
$ws.Range("J17").Value = 436
$ws.Range("K17").Value = 3.669724770642
$ws.Range("L17").Value = 37.386018237082
$ws.Range("M17").Value = 86.776859504132
$ws.Range("N17").Value = -41.902313624678

# Row 18: Burglary
$ws.Range("C18").Value = 47
$ws.Range("D18").Value = 57
$ws.Range("E18").Value = -17.543859649122
$ws.Range("F18").Value = 168
$ws.Range("G18").Value = 166
$ws.Range("H18").Value = 1.204819277108
$ws.Range("I18").Value = 285
$ws.Range("J18").Value = 259
$ws.Range("K18").Value = 10.03861003861
$ws.Range("L18").Value = 41.089108910891
$ws.Range("M18").Value = -37.224669603524
$ws.Range("N18").Value = -89.571899012074

# Row 19: Gr. Larceny
$ws.Range("C19").Value = 100
$ws.Range("D19").Value = 133
$ws.Range("E19").Value = -24.812030075188
$ws.Range("F19").Value = 478
$ws.Range("G19").Value = 519
$ws.Range("H19").Value = -7.899807321772
$ws.Range("I19").Value = 914
$ws.Range("J19").Value = 983
$ws.Range("K19").Value = -7.019328585961
$ws.Range("L19").Value = 71.482176360225
$ws.Range("M19").Value = 28.011204481792
$ws.Range("N19").Value = -16.146788990825

# Row 20: G.L.A.
$ws.Range("C20").Value = 37
$ws.Range("D20").Value = 53
$ws.Range("E20").Value = -30.188679245283
$ws.Range("F20").Value = 129
$ws.Range("G20").Value = 156
$ws.Range("H20").Value = -17.307692307692
$ws.Range("I20").Value = 207
$ws.Range("J20").Value = 253
$ws.Range("K20").Value = -18.181818181818
$ws.Range("L20").Value = 37.086092715231
$ws.Range("M20").Value = -25.270758122743
$ws.Range("N20").Value = -93.899204244031

# Row 21: TOTAL
$ws.Range("C21").Value = 290
$ws.Range("D21").Value = 382
$ws.Range("E21").Value = -24.083769633507
$ws.Range("F21").Value = 1169
$ws.Range("G21").Value = 1303
$ws.Range("H21").Value = -10.283960092095
$ws.Range("I21").Value = 2157
$ws.Range("J21").Value = 2255
$ws.Range("K21").Value = -4.345898004434
$ws.Range("L21").Value = 53.413940256045
$ws.Range("M21").Value = 3.205741626794
$ws.Range("N21").Value = -78.734102336586

# Row 22: Transit
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -33.333333333333
$ws.Range("F22").Value = 12
$ws.Range("H22").Value = -14.285714285714
$ws.Range("I22").Value = 21
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = -16
$ws.Range("L22").Value = 40
$ws.Range("M22").Value = -38.235294117647

# Row 23: Housing
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 12
$ws.Range("E23").Value = -41.666666666666
$ws.Range("F23").Value = 24
$ws.Range("G23").Value = 30
$ws.Range("H23").Value = -20
$ws.Range("I23").Value = 60
$ws.Range("J23").Value = 53
$ws.Range("K23").Value = 13.207547169811
$ws.Range("L23").Value = -6.25
$ws.Range("M23").Value = 76.470588235294

# Row 24: Petit Larceny
$ws.Range("C24").Value = 287
$ws.Range("D24").Value = 274
$ws.Range("E24").Value = 4.744525547445
$ws.Range("F24").Value = 1173
$ws.Range("G24").Value = 1097
$ws.Range("H24").Value = 6.927985414767
$ws.Range("I24").Value = 2141
$ws.Range("J24").Value = 1860
$ws.Range("K24").Value = 15.107526881720
$ws.Range("L24").Value = 36.456341618865
$ws.Range("M24").Value = 35.764109067850

# Row 25: Misd. Assault
$ws.Range("C25").Value = 96
$ws.Range("D25").Value = 95
$ws.Range("E25").Value = 1.052631578947
$ws.Range("F25").Value = 411
$ws.Range("G25").Value = 384
$ws.Range("H25").Value = 7.03125
$ws.Range("I25").Value = 728
$ws.Range("J25").Value = 645
$ws.Range("K25").Value = 12.868217054263
$ws.Range("L25").Value = 35.567970204841
$ws.Range("M25").Value = -10.674846625766

# Row 26: UCR Rape*
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -14.285714285714
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 46
$ws.Range("J26").Value = 40
$ws.Range("K26").Value = 15
$ws.Range("L26").Value = 4.545454545454

# Row 27: Other Sex Crimes
$ws.Range("C27").Value = 9
$ws.Range("D27").Value = 12
$ws.Range("E27").Value = -25
$ws.Range("G27").Value = 49
$ws.Range("H27").Value = -12.244897959183
$ws.Range("I27").Value = 72
$ws.Range("J27").Value = 81
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = 16.129032258064

# Row 28: Shooting Vic.
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 19
$ws.Range("H28").Value = -21.052631578947
$ws.Range("I28").Value = 25
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = -16.666666666666
$ws.Range("L28").Value = 78.571428571428
$ws.Range("M28").Value = 25
$ws.Range("N28").Value = -74.489795918367

# Row 29: Shooting Inc.
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 12
$ws.Range("G29").Value = 16
$ws.Range("H29").Value = -25
$ws.Range("I29").Value = 21
$ws.Range("J29").Value = 23
$ws.Range("K29").Value = -8.695652173913
$ws.Range("L29").Value = 50
$ws.Range("M29").Value = 16.666666666666
$ws.Range("N29").Value = -76.136363636363

# Row 30: Hate Crimes
$ws.Range("D30").Value = 3
$ws.Range("G30").Value = 12
$ws.Range("H30").Value = -83.333333333333
$ws.Range("J30").Value = 18
$ws.Range("K30").Value = -83.333333333333
